# Fill in missing time-tracking numbers for task ID 13 ("Metoda remove")
# across the four sheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "Poraba časa" (Tabela1) : row 14 = ID 13 -----------------------
$ws1 = $wb.Worksheets.Item("Poraba časa")
$ws1.Range("E14").Value = 12
$ws1.Range("F14").Value = 20
$ws1.Range("G14").Value = 0
$ws1.Range("H14").Value = 0
$ws1.Range("I14").Value = 7
$ws1.Range("J14").Value = 0
$ws1.Range("K14").Select()

# --- Sheet "Opravila" (Tabela3) : row 14 = ID 13 --------------------------
$ws2 = $wb.Worksheets.Item("Opravila")
$ws2.Range("C14").Value = 63
$ws2.Range("D14").Value = 29

# --- Sheet "Napake" (Tabela2) : add row for ID 13 -------------------------
$ws3 = $wb.Worksheets.Item("Napake")
$ws3.Range("A4").Value = 13
$ws3.Range("B4").Value = "Metoda remove"
$ws3.Range("C4").Value = 0
$ws3.Range("D4").Value = 0
$ws3.Range("E4").Value = 0
$ws3.Range("F4").Value = 1
$ws3.Range("G4").Value = 0
$ws3.Range("H4").Value = 0
$ws3.Range("I4").Value = 0
$ws3.Range("J4").Value = 0
$ws3.Range("K4").Formula = "=SUM(Tabela2[[#This Row],[Planiranje]:[Analiza]])"
$ws3.Range("A4:B4").Select()

# --- Sheet "Odpravljene napake" (Tabela5) : fill blanks + add row --------
$ws4 = $wb.Worksheets.Item("Odpravljene napake")
$ws4.Range("J2").Value = 0

$ws4.Range("C3").Value = 0
$ws4.Range("D3").Value = 0
$ws4.Range("E3").Value = 0
$ws4.Range("F3").Value = 0
$ws4.Range("H3").Value = 0
$ws4.Range("J3").Value = 0

$ws4.Range("A4").Value = 13
$ws4.Range("B4").Value = "Metoda remove"
$ws4.Range("C4").Value = 0
$ws4.Range("D4").Value = 0
$ws4.Range("E4").Value = 0
$ws4.Range("F4").Value = 0
$ws4.Range("G4").Value = 0
$ws4.Range("H4").Value = 0
$ws4.Range("I4").Value = 1
$ws4.Range("J4").Value = 0
$ws4.Range("K4").Formula = "=SUM(Tabela5[[#This Row],[Planiranje]:[Analiza]])"
$ws4.Range("K4").Select()
